# Auto-generated edit script applying Universalis market-data refresh
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1105.2307
$ws.Range("I12").Value = 2098
$ws.Range("J12").Value = 807.4
$ws.Range("K12").Value = 2098
$ws.Range("L12").Value = 807.4
$ws.Range("M12").Value = -1928
$ws.Range("N12").Value = -1147.4
$ws.Range("H18").Value = 862.25
$ws.Range("I18").Value = 924.5
$ws.Range("K18").Value = 924.5
$ws.Range("M18").Value = -640.5
$ws.Range("H33").Value = 910.3077
$ws.Range("I33").Value = 271.22223
$ws.Range("J33").Value = 2348.25
$ws.Range("K33").Value = 271.22223
$ws.Range("L33").Value = 2348.25
$ws.Range("M33").Value = -42.22223000000002
$ws.Range("N33").Value = -2806.25
$ws.Range("H64").Value = 13248.125
$ws.Range("I64").Value = 9998
$ws.Range("J64").Value = 18665
$ws.Range("K64").Value = 9998
$ws.Range("L64").Value = 18665
$ws.Range("M64").Value = -9750
$ws.Range("N64").Value = -19161
$ws.Range("H67").Value = 13248.125
$ws.Range("I67").Value = 9998
$ws.Range("J67").Value = 18665
$ws.Range("K67").Value = 9998
$ws.Range("L67").Value = 18665
$ws.Range("M67").Value = -9140
$ws.Range("N67").Value = -20381
$ws.Range("H115").Value = 1183.3334
$ws.Range("I115").Value = 1183.3334
$ws.Range("K115").Value = 3550.0002
$ws.Range("M115").Value = -1983.0002
$ws.Range("H127").Value = 700
$ws.Range("I127").Value = 700
$ws.Range("K127").Value = 2100
$ws.Range("M127").Value = 2860
$ws.Range("H132").Value = 1997.1578
$ws.Range("I132").Value = 1929.7333
$ws.Range("K132").Value = 5789.199900000001
$ws.Range("M132").Value = -3259.199900000001
$ws.Range("H138").Value = 13210.2
$ws.Range("J138").Value = 13248.482
$ws.Range("L138").Value = 39745.446
$ws.Range("N138").Value = -50025.446

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3522.1538
$ws.Range("I61").Value = 2112.8572
$ws.Range("K61").Value = 2112.8572
$ws.Range("M61").Value = -1900.8572
$ws.Range("H102").Value = 1499.8334
$ws.Range("I102").Value = 1499.8334
$ws.Range("K102").Value = 1499.8334
$ws.Range("M102").Value = 122.1666
$ws.Range("H131").Value = 90000
$ws.Range("J131").Value = 90000
$ws.Range("L131").Value = 90000
$ws.Range("N131").Value = -100080
$ws.Range("H136").Value = 3522.1538
$ws.Range("I136").Value = 2112.8572
$ws.Range("K136").Value = 6338.571599999999
$ws.Range("M136").Value = -3788.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2959.1904
$ws.Range("I20").Value = 2264.8333
$ws.Range("K20").Value = 2264.8333
$ws.Range("M20").Value = -2017.8333
$ws.Range("H105").Value = 4797.6665
$ws.Range("I105").Value = 3938.8572
$ws.Range("K105").Value = 3938.8572
$ws.Range("M105").Value = -2191.8572
$ws.Range("H134").Value = 1814.2858
$ws.Range("I134").Value = 1450
$ws.Range("K134").Value = 4350
$ws.Range("M134").Value = -1815

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 598.875
$ws.Range("I5").Value = 558.3333
$ws.Range("J5").Value = 623.2
$ws.Range("K5").Value = 558.3333
$ws.Range("L5").Value = 623.2
$ws.Range("M5").Value = -446.3333
$ws.Range("N5").Value = -847.2
$ws.Range("H16").Value = 729.6667
$ws.Range("I16").Value = 729.6667
$ws.Range("K16").Value = 729.6667
$ws.Range("M16").Value = -442.6667
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H93").Value = 20615.625
$ws.Range("I93").Value = 20615.625
$ws.Range("K93").Value = 20615.625
$ws.Range("M93").Value = -18743.625
$ws.Range("H99").Value = 7134.7856
$ws.Range("I99").Value = 5543.5557
$ws.Range("K99").Value = 5543.5557
$ws.Range("M99").Value = -4045.5557
$ws.Range("H113").Value = 729.6667
$ws.Range("I113").Value = 729.6667
$ws.Range("K113").Value = 729.6667
$ws.Range("M113").Value = 1440.3333
$ws.Range("H126").Value = 7134.7856
$ws.Range("I126").Value = 5543.5557
$ws.Range("K126").Value = 16630.6671
$ws.Range("M126").Value = -14160.6671
$ws.Range("H133").Value = 75000
$ws.Range("J133").Value = 75000
$ws.Range("L133").Value = 75000
$ws.Range("N133").Value = -80060
$ws.Range("H134").Value = 3401.6956
$ws.Range("I134").Value = 2697.3125
$ws.Range("K134").Value = 8091.9375
$ws.Range("M134").Value = -5556.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2177.2222
$ws.Range("I122").Value = 2140.5
$ws.Range("K122").Value = 19264.5
$ws.Range("M122").Value = -16814.5
$ws.Range("H131").Value = 1713.7142
$ws.Range("I131").Value = 1319.2
$ws.Range("J131").Value = 1932.8889
$ws.Range("K131").Value = 3957.6
$ws.Range("L131").Value = 5798.6667
$ws.Range("M131").Value = 1082.4
$ws.Range("N131").Value = -15878.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 34718.266
$ws.Range("J24").Value = 34718.266
$ws.Range("L24").Value = 34718.266
$ws.Range("N24").Value = -35064.266
$ws.Range("H70").Value = 16670287
$ws.Range("I70").Value = 33336236
$ws.Range("J70").Value = 4338
$ws.Range("K70").Value = 33336236
$ws.Range("L70").Value = 4338
$ws.Range("M70").Value = -33335966
$ws.Range("N70").Value = -4878
$ws.Range("H73").Value = 16670287
$ws.Range("I73").Value = 33336236
$ws.Range("J73").Value = 4338
$ws.Range("K73").Value = 33336236
$ws.Range("L73").Value = 4338
$ws.Range("M73").Value = -33335300
$ws.Range("N73").Value = -6210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2143.75
$ws.Range("I7").Value = 2075.4
$ws.Range("K7").Value = 2075.4
$ws.Range("M7").Value = -1963.4
$ws.Range("H10").Value = 1669.2858
$ws.Range("I10").Value = 1003
$ws.Range("J10").Value = 1935.8
$ws.Range("K10").Value = 1003
$ws.Range("L10").Value = 1935.8
$ws.Range("M10").Value = -863
$ws.Range("N10").Value = -2215.8
$ws.Range("H22").Value = 1180.375
$ws.Range("J22").Value = 1515.1666
$ws.Range("L22").Value = 1515.1666
$ws.Range("N22").Value = -2105.1666
$ws.Range("H27").Value = 1180.375
$ws.Range("J27").Value = 1515.1666
$ws.Range("L27").Value = 1515.1666
$ws.Range("N27").Value = -1729.1666
$ws.Range("H68").Value = 1500
$ws.Range("J68").Value = 1500
$ws.Range("L68").Value = 1500
$ws.Range("N68").Value = -2998
$ws.Range("H71").Value = 1500
$ws.Range("J71").Value = 1500
$ws.Range("L71").Value = 7500
$ws.Range("N71").Value = -14988
$ws.Range("H116").Value = 289666.66
$ws.Range("J116").Value = 289666.66
$ws.Range("L116").Value = 289666.66
$ws.Range("N116").Value = -298844.66
$ws.Range("H126").Value = 2143.75
$ws.Range("I126").Value = 2075.4
$ws.Range("K126").Value = 6226.200000000001
$ws.Range("M126").Value = -3756.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1954.4546
$ws.Range("I96").Value = 1555.5555
$ws.Range("K96").Value = 1555.5555
$ws.Range("M96").Value = -182.5554999999999
$ws.Range("H122").Value = 1139.8
$ws.Range("I122").Value = 1099.6666
$ws.Range("K122").Value = 3298.9998
$ws.Range("M122").Value = -848.9998000000001
$ws.Range("H126").Value = 2017.1111
$ws.Range("I126").Value = 1450.7142
$ws.Range("J126").Value = 3999.5
$ws.Range("K126").Value = 4352.142599999999
$ws.Range("L126").Value = 11998.5
$ws.Range("M126").Value = -1882.142599999999
$ws.Range("N126").Value = -16938.5
